$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (pushes existing E "information" column to F,
# shifts styles/col widths/trailing empty-format cells along with it).
$ws.Columns("E:E").Insert()

# Header for the newly inserted column.
$ws.Range("E1").Value = "county"

# County values for each center row, keyed by row number.
$counties = @{
    2  = "King County"
    3  = "King County"
    4  = "Snohomish County"
    5  = "Kitsap County"
    6  = "King County"
    7  = "Snohomish County"
    8  = "King County"
    9  = "King County"
    10 = "King County"
    11 = "King County"
    12 = "Pierce County"
    13 = "Snohomish County"
    14 = "Pierce County"
    15 = "Pierce County"
    16 = "King County"
    17 = "King County"
    18 = "King County"
    19 = "King County"
    20 = "King County"
    21 = "King County"
    22 = "King County"
    23 = "King County"
    24 = "King County"
    25 = "King County"
    26 = "Kitsap County"
    27 = "Pierce County"
    28 = "Pierce County"
    29 = "King County"
    30 = "Pierce County"
    31 = "King County"
    32 = "Snohomish County"
    33 = "King County"
    34 = "Pierce County"
    35 = "King County"
    36 = "King County"
    37 = "Pierce County"
    38 = "Snohomish County"
    39 = "Kitsap County"
    40 = "Pierce County"
    41 = "King County"
}

foreach ($row in $counties.Keys) {
    $ws.Cells.Item($row, 5).Value = $counties[$row]
}

# Row 41 gains the same row height as the other data rows (58pt).
$ws.Rows("41").RowHeight = 58
